$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# E56 / E57: previously stored as inline-string text ("543287"/"590024"),
# now stored as numeric values.
$ws.Range("E56").Value = 543287
$ws.Range("E57").Value = 590024

# New row 58
$ws.Range("A58").Value = "26/06/2024 04:44:43"
$ws.Range("B58").Value = 1
$ws.Range("C58").Value = "FACT"
$ws.Range("D58").Value = "Fertilizers And Chemicals Travancore Limited"
$ws.Range("E58").Value = "'590024"
$ws.Range("F58").Value = -0.28
$ws.Range("G58").Value = 996.25
$ws.Range("H58").Value = 551434

# New row 59
$ws.Range("A59").Value = "26/06/2024 04:44:43"
$ws.Range("B59").Value = 2
$ws.Range("C59").Value = "MAXHEALTH"
$ws.Range("D59").Value = "Max Healthcare Institute Ltd"
$ws.Range("E59").Value = "'543220"
$ws.Range("F59").Value = -1.99
$ws.Range("G59").Value = 875.7
$ws.Range("H59").Value = 692127
